$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testTopLinkerUpload")

# Update Top Container Record No. values in column M (rows 6-8): 54555 -> 54556
$ws.Range("M6").Value = 54556
$ws.Range("M7").Value = 54556
$ws.Range("M8").Value = 54556

# Update values in column J (rows 9-10): precreated container ids
$ws.Range("J9").Value = 4
$ws.Range("J10").Value = 4

# Update the active selection on the sheet
$ws.Range("M7:M8").Select()
